{"js": "// Replace each three-digit-division answer string in the document's\n// table cells with its updated value, per the target revision.\nconst replacements = [\n  [\"995\u00f75=199, 0\", \"725\u00f77=103, 4\"],\n  [\"692\u00f75=138, 2\", \"650\u00f76=108, 2\"],\n  [\"438\u00f76=73, 0\", \"928\u00f74=232, 0\"],\n  [\"998\u00f76=166, 2\", \"294\u00f74=73, 2\"],\n  [\"576\u00f72=288, 0\", \"644\u00f77=92, 0\"],\n  [\"267\u00f72=133, 1\", \"784\u00f78=98, 0\"],\n  [\"109\u00f72=54, 1\", \"979\u00f79=108, 7\"],\n  [\"108\u00f76=18, 0\", \"763\u00f77=109, 0\"],\n  [\"406\u00f72=203, 0\", \"732\u00f72=366, 0\"],\n  [\"567\u00f74=141, 3\", \"111\u00f74=27, 3\"],\n  [\"696\u00f77=99, 3\", \"487\u00f77=69, 4\"],\n  [\"557\u00f72=278, 1\", \"190\u00f77=27, 1\"],\n  [\"922\u00f77=131, 5\", \"322\u00f77=46, 0\"],\n  [\"928\u00f77=132, 4\", \"856\u00f74=214, 0\"],\n  [\"842\u00f73=280, 2\", \"801\u00f73=267, 0\"],\n  [\"715\u00f78=89, 3\", \"171\u00f75=34, 1\"],\n  [\"325\u00f74=81, 1\", \"923\u00f77=131, 6\"],\n  [\"129\u00f72=64, 1\", \"242\u00f72=121, 0\"],\n  [\"288\u00f75=57, 3\", \"231\u00f78=28, 7\"],\n  [\"924\u00f78=115, 4\", \"765\u00f77=109, 2\"],\n  [\"483\u00f74=120, 3\", \"537\u00f73=179, 0\"],\n  [\"511\u00f77=73, 0\", \"867\u00f77=123, 6\"],\n  [\"415\u00f76=69, 1\", \"342\u00f78=42, 6\"],\n  [\"447\u00f79=49, 6\", \"847\u00f78=105, 7\"],\n  [\"301\u00f73=100, 1\", \"695\u00f76=115, 5\"],\n];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n    totalReplaced++;\n  }\n  await context.sync();\n}\n\nreturn `replaced: ${totalReplaced}`;\n", "ps1": "# Replace each three-digit-division answer string in the document's\n# table cells with its updated value, per the target revision.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"995\u00f75=199, 0\", \"725\u00f77=103, 4\"),\n    @(\"692\u00f75=138, 2\", \"650\u00f76=108, 2\"),\n    @(\"438\u00f76=73, 0\", \"928\u00f74=232, 0\"),\n    @(\"998\u00f76=166, 2\", \"294\u00f74=73, 2\"),\n    @(\"576\u00f72=288, 0\", \"644\u00f77=92, 0\"),\n    @(\"267\u00f72=133, 1\", \"784\u00f78=98, 0\"),\n    @(\"109\u00f72=54, 1\", \"979\u00f79=108, 7\"),\n    @(\"108\u00f76=18, 0\", \"763\u00f77=109, 0\"),\n    @(\"406\u00f72=203, 0\", \"732\u00f72=366, 0\"),\n    @(\"567\u00f74=141, 3\", \"111\u00f74=27, 3\"),\n    @(\"696\u00f77=99, 3\", \"487\u00f77=69, 4\"),\n    @(\"557\u00f72=278, 1\", \"190\u00f77=27, 1\"),\n    @(\"922\u00f77=131, 5\", \"322\u00f77=46, 0\"),\n    @(\"928\u00f77=132, 4\", \"856\u00f74=214, 0\"),\n    @(\"842\u00f73=280, 2\", \"801\u00f73=267, 0\"),\n    @(\"715\u00f78=89, 3\", \"171\u00f75=34, 1\"),\n    @(\"325\u00f74=81, 1\", \"923\u00f77=131, 6\"),\n    @(\"129\u00f72=64, 1\", \"242\u00f72=121, 0\"),\n    @(\"288\u00f75=57, 3\", \"231\u00f78=28, 7\"),\n    @(\"924\u00f78=115, 4\", \"765\u00f77=109, 2\"),\n    @(\"483\u00f74=120, 3\", \"537\u00f73=179, 0\"),\n    @(\"511\u00f77=73, 0\", \"867\u00f77=123, 6\"),\n    @(\"415\u00f76=69, 1\", \"342\u00f78=42, 6\"),\n    @(\"447\u00f79=49, 6\", \"847\u00f78=105, 7\"),\n    @(\"301\u00f73=100, 1\", \"695\u00f76=115, 5\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $found = $rng.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: not found -> $old\"\n    }\n}\n\nWrite-Output \"done\"\n"}
